# Update landing zone colors
#
# Column H on Sheet1 holds a hex color string per stratigraphic unit.
# Rows 14-20 are the "landing zone" rows (Wolfcamp XY / A Upper / A Lower /
# B Upper / B Lower / C / D) which previously used a ramp of dark reds.
# Re-color them with a new, more distinguishable palette (the old
# "#CC0000" red is kept, just moved down onto the Wolfcamp B Lower row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H14").Value = "#4D8FD1"   # Wolfcamp XY
$ws.Range("H15").Value = "#4685C2"   # Wolfcamp A Upper
$ws.Range("H16").Value = "#10a588"   # Wolfcamp A Lower
$ws.Range("H17").Value = "#a51078"   # Wolfcamp B Upper
$ws.Range("H18").Value = "#CC0000"   # Wolfcamp B Lower
$ws.Range("H19").Value = "#66CC66"   # Wolfcamp C
$ws.Range("H20").Value = "#5EB85E"   # Wolfcamp D

# H19/H20 lose their inherited centered-alignment style once retyped.
$ws.Range("H19:H20").Style = "Normal"

# Column H width gets recalculated for the new content.
$ws.Columns.Item(8).EntireColumn.AutoFit()

# Leave the selection on the last-edited cell, like the author did.
$ws.Range("H16").Select()
